$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column H width ---
$ws.Columns("H").ColumnWidth = 17.14

# --- Row 1: H1 "Chuyen nhuong" (copy format from G1) ---
$ws.Range("H1").Value = "Chuyển nhượng"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)

# --- Row 2: new formula for G2, new cell H2 ---
$ws.Range("G2").Formula = "=D5+D6+D7+D8+D12+D14+D16+D18+H2+D19"
$ws.Range("H2").Value = -600000
$ws.Range("G1").Copy()
$ws.Range("H2").PasteSpecial(-4122)

# --- Row 3: G3 formula + style change to match F3 (style 3), new cell H3 ---
$ws.Range("G3").Formula = "=D9+D10+D11+D13+D15+D17+H3"
$ws.Range("F3").Copy()
$ws.Range("G3").PasteSpecial(-4122)
$ws.Range("G3").Formula = "=D9+D10+D11+D13+D15+D17+H3"
$ws.Range("H3").Formula = "=-H2"
$ws.Range("G1").Copy()
$ws.Range("H3").PasteSpecial(-4122)
$ws.Range("H3").Formula = "=-H2"

# --- Row 4: G4 style change to match G1/G2 (style 2) ---
$ws.Range("G1").Copy()
$ws.Range("G4").PasteSpecial(-4122)
$ws.Range("G4").Formula = "=G3+G2"

# --- Row 17 ---
$ws.Range("A12").Copy()
$ws.Range("A17").PasteSpecial(-4122)
$ws.Range("A17").Value = 11

$ws.Range("A12").Copy()
$ws.Range("B17").PasteSpecial(-4122)
$ws.Range("B17").Value = "Linh kiện "

$ws.Range("C12").Copy()
$ws.Range("C17").PasteSpecial(-4122)
$ws.Range("C17").Value = "Ốc"

$ws.Range("D12").Copy()
$ws.Range("D17").PasteSpecial(-4122)
$ws.Range("D17").Value = 12000

$ws.Range("E12").Copy()
$ws.Range("E17").PasteSpecial(-4122)
$ws.Range("E17").Value = "Khánh"

$ws.Range("F12").Copy()
$ws.Range("F17").PasteSpecial(-4122)
$ws.Range("F17").Value = "lần 4"

# --- Row 18 ---
$ws.Range("A12").Copy()
$ws.Range("A18").PasteSpecial(-4122)

$ws.Range("A12").Copy()
$ws.Range("B18").PasteSpecial(-4122)

$ws.Range("C12").Copy()
$ws.Range("C18").PasteSpecial(-4122)
$ws.Range("C18").Value = "lần 5"

$ws.Range("D12").Copy()
$ws.Range("D18").PasteSpecial(-4122)
$ws.Range("D18").Value = 14000

$ws.Range("E12").Copy()
$ws.Range("E18").PasteSpecial(-4122)
$ws.Range("E18").Value = "Long"

$ws.Range("F12").Copy()
$ws.Range("F18").PasteSpecial(-4122)
$ws.Range("F18").Value = "Keo dán"

# --- Row 19 ---
$ws.Range("A12").Copy()
$ws.Range("A19").PasteSpecial(-4122)
$ws.Range("A19").HorizontalAlignment = 1
$ws.Range("A19").Value = 12

$ws.Range("C12").Copy()
$ws.Range("B19").PasteSpecial(-4122)
$ws.Range("B19").Value = "Module cầu H"

$ws.Range("C12").Copy()
$ws.Range("C19").PasteSpecial(-4122)
$ws.Range("C19").Value = "1 cái"

$ws.Range("D12").Copy()
$ws.Range("D19").PasteSpecial(-4122)
$ws.Range("D19").Value = 150000

$ws.Range("E12").Copy()
$ws.Range("E19").PasteSpecial(-4122)
$ws.Range("E19").Value = "Long"

$ws.Range("F12").Copy()
$ws.Range("F19").PasteSpecial(-4122)
$ws.Range("F19").Value = "ngày 24/03/2016"

# --- Merge new cells ---
$ws.Range("A17:A18").Merge()
$ws.Range("B17:B18").Merge()

# --- Selection ---
$ws.Range("I5").Select()
